# Weapons Adjustments + 60fps cap
# + All Weapons Range Reduction
# + Remade Blasters DMG
# + Linear Paint Droplets (For Blasters)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4 (Blasters - default) ---
$ws.Range("E4").Value = 10
$ws.Range("K4").Value = 10
$ws.Range("P4").Value = 1

# --- Row 5 (Blasters - express) ---
$ws.Range("E5").Value = 12
$ws.Range("K5").Value = 10
$ws.Range("P5").Value = 1

# --- Row 6 (Blasters - heavy) ---
$ws.Range("E6").Value = 6.67
$ws.Range("K6").Value = 12
$ws.Range("P6").Value = 1

# --- Row 7 (Blasters - long) ---
$ws.Range("E7").Value = 8
$ws.Range("K7").Value = 16
$ws.Range("P7").Value = 2
$ws.Range("Q7").Value = 1.75

# --- Row 8 (Blasters - rng / shotgun style) ---
$ws.Range("C8").Value = 38
$ws.Range("E8").Value = 12
$ws.Range("K8").Value = 7
$ws.Range("O8").Value = 4
$ws.Range("P8").Value = 1

# --- Row 9 (Shotters - burst_short) ---
$ws.Range("K9").Value = 13.5

# --- Row 10 (Shotters - burst_long) ---
$ws.Range("K10").Value = 16.5

# --- Row 11 (Miniguns - default) ---
$ws.Range("C11").Value = 13
$ws.Range("D11").Value = 7
$ws.Range("E11").Value = 2.25
$ws.Range("K11").Value = 11

# --- Row 12 (Miniguns - express) ---
$ws.Range("E12").Value = 3
$ws.Range("K12").Value = 7.5
$ws.Range("Q12").Value = 2.5

# --- Row 13 (Miniguns - heavy) ---
$ws.Range("K13").Value = 17

# --- Row 14 (Blasters [fallof or] - long) ---
$ws.Range("K14").Value = 16
$ws.Range("O14").Value = 4

# --- Row 15 (Blasters express) ---
$ws.Range("K15").Value = 16
$ws.Range("O15").Value = 2.5

# --- Row 16 (Blasters heavy) ---
$ws.Range("O16").Value = 4.25

# --- Row 17 (Miniguns 70-35 default) ---
$ws.Range("K17").Value = 17
$ws.Range("P17").Value = 1

# --- Row 18 (express) ---
$ws.Range("K18").Value = 13
$ws.Range("P18").Value = 1

# --- Row 19 (heavy) ---
$ws.Range("K19").Value = 20

# --- Recreate shared/fill-down formulas so the dependent ranges
#     pick up consistent shared formula groups, matching how Excel
#     records a fill-down edit across these weapon rows. ---

# R4:R19  => O*D + P*Q*D
$ws.Range("R4:R19").Formula = "=O4*D4+P4*Q4*D4"

# S4:S19  => E*R
$ws.Range("S4:S19").Formula = "=E4*R4"

# G11:G19 => C*E*D
$ws.Range("G11:G19").Formula = "=C11*E11*D11"

# H11:H19 => 100/C
$ws.Range("H11:H19").Formula = "=100/C11"

# J11:J19 => 100/I
$ws.Range("J11:J19").Formula = "=100/I11"

# F12:F17 => 1/E (fill down from F12, leaving F11 standalone)
$ws.Range("F12:F17").Formula = "=1/E12"

# --- Update active cell / selection to match the saved cursor position ---
$ws.Range("E12").Select() | Out-Null
